$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1B Summary")

# Clear the value in B6 (was a literal 0, should be empty while keeping style s="16")
$ws.Range("B6").ClearContents()
